$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row (row 1) from the "_old"/"_new" suffix convention
#    to the "_FV2310"/"_FV2404" suffix convention (column "diff" is unchanged).
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2310"
}

# Column 11 ("diff") keeps its value as-is.

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2404"
}

# 2) Turn the A1:U58 range into an Excel Table ("Table1") now that the
#    header cells carry their final names.
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# 3) Freeze the header row (split/freeze pane above row 2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
